# Doing Updates for Financials
# Update yearly financial figures on the TPIC sheet (most-recent-period column D,
# plus one figure in column I for Capital Expenditures).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TPIC")

# Income Statement
$ws.Range("D8").Value  = 955200    # Total Revenue
$ws.Range("D9").Value  = 844700    # Cost of Revenue
$ws.Range("D10").Value = 110500    # Gross Profit
$ws.Range("D17").Value = 885100    # Total Operating Expenses
$ws.Range("D18").Value = 70100     # Operating Income or Loss
$ws.Range("D21").Value = 87800     # Earnings Before Interest And Taxes
$ws.Range("D23").Value = 54500     # Income Before Tax
$ws.Range("D24").Value = 15800     # Income Tax Expense
$ws.Range("D26").Value = 38700     # Income After Tax
$ws.Range("D27").Value = 38700     # Net Income From Continuing Ops
$ws.Range("D33").Value = -35600    # Net Income
$ws.Range("D35").Value = -35600    # Net Income Applicable To Common Shares

# Balance Sheet
$ws.Range("D43").Value = 348800    # Net Receivables
$ws.Range("D45").Value = 58900     # Other Current Assets
$ws.Range("D48").Value = 247000    # Property Plant and Equipment
$ws.Range("D52").Value = 25600     # Other Assets
$ws.Range("D54").Value = 545700    # Total Assets
$ws.Range("D62").Value = 3400      # Other Liabilities
$ws.Range("D66").Value = 325200    # Total Liabilities
$ws.Range("D76").Value = 220600    # Total Stockholder Equity

# Cash Flow Statement
$ws.Range("D81").Value = -35600    # Net Income
$ws.Range("I91").Value = -4300     # Capital Expenditures
